$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing attendee row, shifting it down to row 3.
$ws.Rows.Item(2).Insert()

# Fill the newly inserted row 2 with the latest check-in.
$ws.Cells.Item(2, 1).Value = "HE170769"
$ws.Cells.Item(2, 2).Value = "2023-09-08 17:01:29"

# Update the timestamp for the attendee that was pushed down to row 3.
$ws.Cells.Item(3, 2).Value = "2023-09-08 17:01:50"
